$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (New York -- New York): fill in previously-empty scraped data ---

# B4: Date Published (serial date 44034 -> 2020-07-22), same date format as other rows
$ws.Range("B4").Value = 44034
$ws.Range("B4").NumberFormat = "YYYY-MM-DD"

# C4 / D4: these come through the scrape as text (numeric-looking strings),
# matching the pattern already present elsewhere in this sheet (e.g. C8/D8).
$ws.Range("C4").Value = "'219128"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = "'18803"
$ws.Range("D4").Style = "Normal"

# E4 - H4: numeric counts/percentages
$ws.Range("E4").Value = 33790
$ws.Range("F4").Value = 5239
$ws.Range("G4").Value = 30.07
$ws.Range("H4").Value = 30.43

# J4: Pct Includes Hispanic Black -> TRUE
$ws.Range("J4").Value = $true

# K4 / L4: numeric counts
$ws.Range("K4").Value = 112360
$ws.Range("L4").Value = 17217

# O4: Status code -> Success!
$ws.Range("O4").Value = "Success!"

# --- Row 41 (Iowa): updated Count Cases Black/AA ---
$ws.Range("E41").Value = 3288
